$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 '25.926.00'
Set-TextCell 2 5 '  -0.22%  '

# Row 3
Set-TextCell 3 4 '1.638.82'
Set-TextCell 3 5 '  -0.29%  '

# Row 4
Set-TextCell 4 4 '1.003'
Set-TextCell 4 5 '  +0.17%  '

# Row 5
Set-TextCell 5 4 '214.71'
Set-TextCell 5 5 '  -0.39%  '

# Row 6
Set-TextCell 6 4 '0.5055'
Set-TextCell 6 5 '  -0.13%  '

# Row 7
Set-TextCell 7 4 '1.002'
Set-TextCell 7 5 '  +0.07%  '

# Row 8
Set-TextCell 8 5 '  -0.99%  '

# Row 9
Set-TextCell 9 4 '0.06366'
Set-TextCell 9 5 '  -0.85%  '

# Row 10
Set-TextCell 10 5 '  -0.91%  '

# Row 11
Set-TextCell 11 4 '0.07754'
Set-TextCell 11 5 '  -0.17%  '

# Row 12
Set-TextCell 12 4 '4.275'
Set-TextCell 12 5 '  +0.03%  '

# Row 13
Set-TextCell 13 4 '1.644.35'
Set-TextCell 13 5 '  +0.00%  '

# Row 14
Set-TextCell 14 4 '0.5440'
Set-TextCell 14 5 '  -0.54%  '

# Row 15
Set-TextCell 15 4 '0.0₅7803'
Set-TextCell 15 5 '  -1.69%  '

# Row 16
Set-TextCell 16 4 '64.19'
Set-TextCell 16 5 '  -0.40%  '

# Row 17
Set-TextCell 17 4 '25.957.10'
Set-TextCell 17 5 '  -0.18%  '

# Row 18
Set-TextCell 18 4 '1.002'
Set-TextCell 18 5 '  +0.07%  '

# Row 19
Set-TextCell 19 4 '197.11'
Set-TextCell 19 5 '  -2.68%  '

# Row 20
Set-TextCell 20 4 '4.449'
Set-TextCell 20 5 '  +1.36%  '

# Row 21
Set-TextCell 21 4 '9.928'
Set-TextCell 21 5 '  +0.14%  '

# Row 22
Set-TextCell 22 4 '6.022'
Set-TextCell 22 5 '  +0.52%  '

# Row 23
Set-TextCell 23 4 '1.004'
Set-TextCell 23 5 '  +0.20%  '

# Row 24
Set-TextCell 24 4 '1.887'
Set-TextCell 24 5 '  +0.17%  '

# Row 25
Set-TextCell 25 4 '140.63'
Set-TextCell 25 5 '  -0.36%  '

# Row 26
Set-TextCell 26 4 '0.1177'
Set-TextCell 26 5 '  +3.21%  '

# Row 27
Set-TextCell 27 4 '6.867'
Set-TextCell 27 5 '  +0.66%  '

# Row 28
Set-TextCell 28 4 '15.70'
Set-TextCell 28 5 '  -0.18%  '

# Row 29
Set-TextCell 29 5 '  -0.61%  '

# Row 30
Set-TextCell 30 4 '0.04926'
Set-TextCell 30 5 '  +0.10%  '

# Row 31
Set-TextCell 31 4 '3.250'
Set-TextCell 31 5 '  -0.74%  '

# Row 32
Set-TextCell 32 4 '3.176'
Set-TextCell 32 5 '  -1.30%  '

# Row 33
Set-TextCell 33 5 '  -0.60%  '

# Row 34
Set-TextCell 34 4 '2.363'
Set-TextCell 34 5 '  +0.03%  '

# Row 35
Set-TextCell 35 4 '0.8928'
Set-TextCell 35 5 '  -0.14%  '

# Row 36
Set-TextCell 36 4 '2.578'
Set-TextCell 36 5 '  -1.82%  '

# Row 37
Set-TextCell 37 4 '1.132.30'
Set-TextCell 37 5 '  -1.88%  '

# Row 38
Set-TextCell 38 4 '0.5429'
Set-TextCell 38 5 '  -3.00%  '

# Row 39
Set-TextCell 39 4 '0.01556'
Set-TextCell 39 5 '  -0.75%  '

# Row 40
Set-TextCell 40 4 '2.551'
Set-TextCell 40 5 '  -0.19%  '

# Row 41
Set-TextCell 41 4 '1.002'
Set-TextCell 41 5 '  +0.16%  '

# Row 42
Set-TextCell 42 2 'BabyDogeCoin'
Set-TextCell 42 3 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell 42 4 '0.0₈129'
Set-TextCell 42 5 '  +9.21%  '

# Row 43
Set-TextCell 43 2 'FraxShare'
Set-TextCell 43 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 43 4 '5.572'
Set-TextCell 43 5 '  -2.49%  '

# Row 44
Set-TextCell 44 2 'TrustWalletToken'
Set-TextCell 44 3 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell 44 4 '0.8171'
Set-TextCell 44 5 '  +1.43%  '

# Row 45
Set-TextCell 45 4 '99.33'
Set-TextCell 45 5 '  -0.48%  '

# Row 46
Set-TextCell 46 4 '1.777.38'
Set-TextCell 46 5 '  -0.19%  '

# Row 47
Set-TextCell 47 4 '0.4535'
Set-TextCell 47 5 '  +0.49%  '

# Row 48
Set-TextCell 48 4 '1.003'
Set-TextCell 48 5 '  -0.20%  '

# Row 49
Set-TextCell 49 4 '54.70'
Set-TextCell 49 5 '  -0.15%  '

# Row 50
Set-TextCell 50 4 '0.05073'
Set-TextCell 50 5 '  +0.54%  '

# Row 51
Set-TextCell 51 5 '  +0.31%  '
